$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value = 160
$ws.Cells.Item(7, 6).Value = 941
$ws.Cells.Item(9, 6).Value = 183
$ws.Cells.Item(11, 6).Value = 80
$ws.Cells.Item(15, 6).Value = 491
$ws.Cells.Item(16, 6).Value = 1290
$ws.Cells.Item(18, 6).Value = 423
$ws.Cells.Item(19, 6).Value = 1083
$ws.Cells.Item(20, 6).Value = 2793
$ws.Cells.Item(22, 6).Value = 646
$ws.Cells.Item(23, 6).Value = 166
$ws.Cells.Item(24, 6).Value = 1241
$ws.Cells.Item(28, 6).Value = 826
$ws.Cells.Item(29, 6).Value = 12
$ws.Cells.Item(31, 6).Value = 1318

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 511

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 718

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 718
$ws.Cells.Item(7, 6).Value = 511
$ws.Cells.Item(8, 6).Value = 511
$ws.Cells.Item(13, 6).Value = 160
$ws.Cells.Item(15, 6).Value = 941
$ws.Cells.Item(17, 6).Value = 183
$ws.Cells.Item(23, 6).Value = 80
$ws.Cells.Item(28, 6).Value = 491
$ws.Cells.Item(29, 6).Value = 1290
$ws.Cells.Item(31, 6).Value = 423
$ws.Cells.Item(32, 6).Value = 1083
$ws.Cells.Item(33, 6).Value = 2793
$ws.Cells.Item(35, 6).Value = 646
$ws.Cells.Item(36, 6).Value = 166
$ws.Cells.Item(37, 6).Value = 1241
$ws.Cells.Item(43, 6).Value = 826
$ws.Cells.Item(44, 6).Value = 12
$ws.Cells.Item(46, 6).Value = 1318

$wb.Save()
